# Insert a "Docentes responsáveis:" block (label + 3 names) right after the
# "Objectives:" row (row 11), pushing every following row down by 4.
#
# Final layout for the new rows:
#   row 12: A = "Docentes responsáveis:"
#   row 13: B/C = "7459752 - Maria Ismenia Sodero Toledo Faria"
#   row 14: B/C = "2166002 - Sandra Giacomin Schneider"
#   row 15: B/C = "1922320 - Sebastiao Ribeiro"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 12 (old row 12 "Programa resumido:" and
# everything after it shifts down to row 16 onward).
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(12).Insert()
}

# The Insert() above stamps every new row with a copy of column A's
# formatting (from the row that used to be at 12). Rows 13-15 should have no
# value in column A at all, so clear those cells completely.
$ws.Range("A13:A15").Clear()

# Row 12: the new section label, bold (same style as other column-A labels).
$ws.Range("A12").Value = "Docentes responsáveis:"

# Rows 13-15: teacher names duplicated into columns B and C, matching the
# pattern used by every other two-column content row in the sheet.
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

# Column B on the inserted rows picked up the bold "column A" font from the
# Insert() copy (column C already inherited the right style). Re-apply the
# normal wrap-text style used by every other column-B content cell by
# copying the format from a still-correctly-styled column-B cell.
$ws.Range("B16").Copy()
$ws.Range("B13:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
